# Applies the "point de situation avec la prof" edits to the
# SignUp plan-test document:
#   1. Merge the TC_U101 / TC_U102 / TC_U103 bullet points into a single
#      TC_U101 bullet that documents createUser(data) for "email already
#      exists" validation.
#   2. Rework the TC_U104 bullet into "TC_U102" (same createUser(data) flow,
#      "already used" branch).
#   3. Refresh the tracking table: rows for TC_U101 / TC_U102 get the new
#      description / expected-result text, and the now-redundant TC_U103 /
#      TC_U104 rows are deleted outright.
#   4. Drop the stale lastRenderedPageBreak cached before the
#      "Tableau de suivi des tests" heading (content reflowed).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Bullet list: merge TC_U101 + TC_U102 + TC_U103 paragraphs into one.
# ---------------------------------------------------------------------
$oldBullet1 = "TC_U101 : validatePassword(password) retourne vrai pour un mot de passe conforme." + [char]13 + `
              "TC_U102 : hashPassword(password) retourne une chaîne hachée non vide." + [char]13 + `
              "TC_U103 : createUser(username, password) crée un utilisateur si le nom n’existe pas déjà."
$newBullet1 = "TC_U101 : createUser(data) crée un utilisateur si l’email n’existe pas déjà."

$d.Content.Find.Execute($oldBullet1, $false, $false, $false, $false, $false, $true, 1, $false, $newBullet1, 2) | Out-Null

# Only "TC_U101 :" stays bold; the rest of the (merged) sentence is plain.
$prefixRng = $d.Content
$prefixRng.Find.Execute("TC_U101 :", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fullRng = $d.Content
$fullRng.Find.Execute($newBullet1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$suffixRng = $d.Range($prefixRng.End, $fullRng.End)
$suffixRng.Font.Bold = 0
$suffixRng.Font.BoldBi = 0

# ---------------------------------------------------------------------
# 2) TC_U104 bullet -> becomes the "TC_U102" createUser(data) / already
#    used branch (paragraph itself is left standalone, only its text
#    changes).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TC_U104 :", $false, $false, $false, $false, $false, $true, 1, $false, "TC_U102 :", 2) | Out-Null
$d.Content.Find.Execute( `
    "createUser(username, password) retourne une erreur si le nom est déjà utilisé.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "createUser(data) retourne une erreur si l’email est déjà utilisé.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Tracking table refresh.
# ---------------------------------------------------------------------
# Row TC_U101: description + expected result.
$d.Content.Find.Execute( `
    "Vérifier la validation d’un mot de passe conforme", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Crée un utilisateur avec un email qui n’existe pas déjà.", 2) | Out-Null

$d.Content.Find.Execute( `
    "La fonction retourne true", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "L’utilisateur est créé", 2) | Out-Null

# Row TC_U102: description + expected result.
$d.Content.Find.Execute( `
    "Vérifier le hachage du mot de passe", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Crée un utilisateur avec un email qui existe déjà.", 2) | Out-Null

$d.Content.Find.Execute( `
    "Le hachage est généré correctement", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Retourne une erreur", 2) | Out-Null

# TC_U103 / TC_U104 rows are gone now (merged into the bullets above /
# no longer relevant) -- drop them from the table, bottom row first so
# indices stay valid.
$tbl = $d.Tables(1)
$tbl.Rows(5).Delete()
$tbl.Rows(4).Delete()

# ---------------------------------------------------------------------
# 4) The heading's cached page-break marker is stale after the reflow;
#    re-typing the run through Find/Replace drops it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Tableau de suivi des tests", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Tableau de suivi des tests", 2) | Out-Null
